$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip thousands-separator commas and the trailing "hp" unit text from the
# "Power Output" column (D2:D60), leaving the bare number as text with a
# trailing space (e.g. "1,200 hp" -> "1200 "). Only touch rows that already
# have a value in column D.
for ($r = 2; $r -le 60; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $text = [string]$val
        $clean = $text.Replace(",", "")
        $clean = $clean -replace "\s*hp\s*$", " "
        $cell.NumberFormat = "@"
        $cell.Value = $clean
        $cell.NumberFormat = "0.00"
    }
}

# D61 becomes a true numeric value (3000) instead of text.
$d61 = $ws.Cells.Item(61, 4)
$d61.Value = 3000
$d61.NumberFormat = "0.00"

# Header D1 also gets the new number format applied.
$ws.Cells.Item(1, 4).NumberFormat = "0.00"

# Set the new column width/format for column D (applies the "0.00" style to
# the whole column, including any currently-empty cells).
$ws.Columns.Item(4).ColumnWidth = 9.140625
$ws.Columns.Item(4).NumberFormat = "0.00"

# Update the active selection to I6
$ws.Range("I6").Select()
